$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @(
  "48-32=",
  "88-26=",
  "54-50=",
  "20+36=",
  "1+49=",
  "50+46=",
  "81-77=",
  "51+22=",
  "99-79=",
  "6+1=",
  "42-30=",
  "21+27=",
  "8-5=",
  "80-74=",
  "83+3=",
  "39-34=",
  "78-17=",
  "50+35=",
  "78-38=",
  "62+16=",
  "7+48=",
  "36+52=",
  "22+64=",
  "90-26=",
  "1+96=",
  "40+14=",
  "60-21=",
  "13+66=",
  "58-14=",
  "15-4=",
  "76+17=",
  "89-14=",
  "39+56=",
  "47-8=",
  "61-56=",
  "83-54=",
  "1+55=",
  "93-63=",
  "90-18=",
  "15+58=",
  "41+19=",
  "30-23=",
  "77-56=",
  "7-1=",
  "64+9=",
  "24-5=",
  "37+12=",
  "81-54=",
  "72+7=",
  "68-32=",
  "45-28=",
  "1+55=",
  "12+35=",
  "84-29=",
  "33+60=",
  "45-1=",
  "75-54=",
  "57+8=",
  "78-72=",
  "17+82=",
  "59+19=",
  "61-29=",
  "0+72=",
  "36+52=",
  "18+27=",
  "11+44=",
  "11+54=",
  "1+66=",
  "17-12=",
  "22+34=",
  "50-37=",
  "58+30=",
  "68+28=",
  "31+31=",
  "8+61=",
  "68-64=",
  "76-22=",
  "68-17=",
  "31-10=",
  "59+8=",
  "46-45=",
  "22-14=",
  "65+18=",
  "38+32=",
  "36-36=",
  "3+73=",
  "13-9=",
  "13+80=",
  "81-11=",
  "11+55=",
  "58-39=",
  "85-9=",
  "30+32=",
  "22+17=",
  "57-10=",
  "82-25=",
  "20+57=",
  "87+11=",
  "48-12=",
  "26+53="
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output "done: updated $idx cells"
